# Applies the data-update edit described by the diff:
#  - refresh the numeric values in every results table (rows 2-8, 12-16,
#    20-26, 38-44) to their new magnitudes
#  - rescale row 50 (E50/F50)
#  - drop the now-unused "6k"-"13k" columns (G:N) from the summary table
#    in rows 48-50 (J48/K48/O48 keep their red-font style but become blank)
#  - update the active selection to L23

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data-table values ---
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = 20
$ws.Range("E2").Value = 30
$ws.Range("F2").Value = 40
$ws.Range("G2").Value = 50
$ws.Range("H2").Value = 60
$ws.Range("I2").Value = 70
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 90
$ws.Range("L2").Value = 100
$ws.Range("C3").Value = 5.8760000000000003
$ws.Range("D3").Value = 8.9429999999999996
$ws.Range("E3").Value = 12.545999999999999
$ws.Range("F3").Value = 15.127000000000001
$ws.Range("G3").Value = 18.576000000000001
$ws.Range("H3").Value = 22.131
$ws.Range("I3").Value = 25.213999999999999
$ws.Range("J3").Value = 28.824999999999999
$ws.Range("K3").Value = 32.673000000000002
$ws.Range("L3").Value = 15.153
$ws.Range("C4").Value = 7.2629999999999999
$ws.Range("D4").Value = 11.763999999999999
$ws.Range("E4").Value = 16.736000000000001
$ws.Range("F4").Value = 21.902999999999999
$ws.Range("G4").Value = 26.152999999999999
$ws.Range("H4").Value = 30.719000000000001
$ws.Range("I4").Value = 34.545999999999999
$ws.Range("J4").Value = 39.457999999999998
$ws.Range("K4").Value = 44.113
$ws.Range("L4").Value = 48.902999999999999
$ws.Range("C5").Value = 4.4530000000000003
$ws.Range("D5").Value = 7.726
$ws.Range("E5").Value = 10.273
$ws.Range("F5").Value = 12.561
$ws.Range("G5").Value = 15.003
$ws.Range("H5").Value = 17.526
$ws.Range("I5").Value = 20.126000000000001
$ws.Range("J5").Value = 22.812999999999999
$ws.Range("K5").Value = 25.135999999999999
$ws.Range("L5").Value = 27.684000000000001
$ws.Range("C6").Value = 16.3
$ws.Range("D6").Value = 32.6
$ws.Range("E6").Value = 48.9
$ws.Range("F6").Value = 65.2
$ws.Range("G6").Value = 81.5
$ws.Range("H6").Value = 97.8
$ws.Range("I6").Value = 114.1
$ws.Range("J6").Value = 130.4
$ws.Range("K6").Value = 146.69999999999999
$ws.Range("L6").Value = 163
$ws.Range("C7").Value = 34
$ws.Range("D7").Value = 68
$ws.Range("E7").Value = 102
$ws.Range("F7").Value = 136
$ws.Range("G7").Value = 170
$ws.Range("H7").Value = 204
$ws.Range("I7").Value = 238
$ws.Range("J7").Value = 272
$ws.Range("K7").Value = 306
$ws.Range("L7").Value = 340
$ws.Range("C8").Value = 22.8
$ws.Range("D8").Value = 45.6
$ws.Range("E8").Value = 68.400000000000006
$ws.Range("F8").Value = 91.2
$ws.Range("G8").Value = 114
$ws.Range("H8").Value = 136.80000000000001
$ws.Range("I8").Value = 159.6
$ws.Range("J8").Value = 182.4
$ws.Range("K8").Value = 205.2
$ws.Range("L8").Value = 228
$ws.Range("C12").Value = 50
$ws.Range("D12").Value = 100
$ws.Range("E12").Value = 150
$ws.Range("F12").Value = 200
$ws.Range("G12").Value = 250
$ws.Range("H12").Value = 300
$ws.Range("I12").Value = 350
$ws.Range("J12").Value = 400
$ws.Range("K12").Value = 450
$ws.Range("L12").Value = 500
$ws.Range("C13").Value = 12.147
$ws.Range("D13").Value = 22.974
$ws.Range("E13").Value = 34.426000000000002
$ws.Range("F13").Value = 45.731999999999999
$ws.Range("G13").Value = 56.137999999999998
$ws.Range("H13").Value = 68.546000000000006
$ws.Range("I13").Value = 79.575000000000003
$ws.Range("J13").Value = 90.016000000000005
$ws.Range("K13").Value = 102.24299999999999
$ws.Range("L13").Value = 113.917
$ws.Range("C14").Value = 14.500999999999999
$ws.Range("D14").Value = 30.390999999999998
$ws.Range("E14").Value = 45.581000000000003
$ws.Range("F14").Value = 60.871000000000002
$ws.Range("G14").Value = 76.212999999999994
$ws.Range("H14").Value = 90.986000000000004
$ws.Range("I14").Value = 105.843
$ws.Range("J14").Value = 119.91500000000001
$ws.Range("K14").Value = 135.137
$ws.Range("L14").Value = 150.136
$ws.Range("C15").Value = 81.5
$ws.Range("D15").Value = 163
$ws.Range("E15").Value = 244.5
$ws.Range("F15").Value = 326
$ws.Range("G15").Value = 407.5
$ws.Range("H15").Value = 489
$ws.Range("I15").Value = 570.5
$ws.Range("J15").Value = 652
$ws.Range("K15").Value = 733.5
$ws.Range("L15").Value = 815
$ws.Range("C16").Value = 114
$ws.Range("D16").Value = 228
$ws.Range("E16").Value = 342
$ws.Range("F16").Value = 456
$ws.Range("G16").Value = 570
$ws.Range("H16").Value = 684
$ws.Range("I16").Value = 798
$ws.Range("J16").Value = 912
$ws.Range("K16").Value = 1026
$ws.Range("L16").Value = 1140
$ws.Range("C20").Value = 50
$ws.Range("D20").Value = 100
$ws.Range("E20").Value = 150
$ws.Range("F20").Value = 200
$ws.Range("G20").Value = 250
$ws.Range("H20").Value = 300
$ws.Range("I20").Value = 350
$ws.Range("J20").Value = 400
$ws.Range("K20").Value = 450
$ws.Range("L20").Value = 500
$ws.Range("C21").Value = 20.268000000000001
$ws.Range("D21").Value = 34.683999999999997
$ws.Range("E21").Value = 47.390999999999998
$ws.Range("F21").Value = 60.527999999999999
$ws.Range("G21").Value = 73.081000000000003
$ws.Range("H21").Value = 84.816999999999993
$ws.Range("I21").Value = 96.861000000000004
$ws.Range("J21").Value = 109.73699999999999
$ws.Range("K21").Value = 121.852
$ws.Range("L21").Value = 134.214
$ws.Range("C22").Value = 25.716999999999999
$ws.Range("D22").Value = 43.594000000000001
$ws.Range("E22").Value = 61.856999999999999
$ws.Range("F22").Value = 77.623000000000005
$ws.Range("G22").Value = 97.070999999999998
$ws.Range("H22").Value = 114.846
$ws.Range("I22").Value = 131.81700000000001
$ws.Range("J22").Value = 147.46600000000001
$ws.Range("K22").Value = 163.625
$ws.Range("L22").Value = 180.983
$ws.Range("C23").Value = 14.987
$ws.Range("D23").Value = 18.917000000000002
$ws.Range("E23").Value = 22.966000000000001
$ws.Range("F23").Value = 27.071000000000002
$ws.Range("G23").Value = 30.911999999999999
$ws.Range("H23").Value = 34.923000000000002
$ws.Range("I23").Value = 38.845999999999997
$ws.Range("J23").Value = 42.566000000000003
$ws.Range("K23").Value = 46.594000000000001
$ws.Range("L23").Value = 50.292999999999999
$ws.Range("C24").Value = 81.5
$ws.Range("D24").Value = 163
$ws.Range("E24").Value = 244.5
$ws.Range("F24").Value = 326
$ws.Range("G24").Value = 407.5
$ws.Range("H24").Value = 489
$ws.Range("I24").Value = 570.5
$ws.Range("J24").Value = 652
$ws.Range("K24").Value = 733.5
$ws.Range("L24").Value = 815
$ws.Range("C25").Value = 170
$ws.Range("D25").Value = 340
$ws.Range("E25").Value = 510
$ws.Range("F25").Value = 680
$ws.Range("G25").Value = 850
$ws.Range("H25").Value = 1020
$ws.Range("I25").Value = 1190
$ws.Range("J25").Value = 1360
$ws.Range("K25").Value = 1530
$ws.Range("L25").Value = 1700
$ws.Range("C26").Value = 114
$ws.Range("D26").Value = 228
$ws.Range("E26").Value = 342
$ws.Range("F26").Value = 456
$ws.Range("G26").Value = 570
$ws.Range("H26").Value = 684
$ws.Range("I26").Value = 798
$ws.Range("J26").Value = 912
$ws.Range("K26").Value = 1026
$ws.Range("L26").Value = 1140
$ws.Range("C38").Value = 30
$ws.Range("D38").Value = 60
$ws.Range("E38").Value = 90
$ws.Range("F38").Value = 120
$ws.Range("G38").Value = 150
$ws.Range("H38").Value = 180
$ws.Range("I38").Value = 210
$ws.Range("J38").Value = 240
$ws.Range("K38").Value = 270
$ws.Range("L38").Value = 300
$ws.Range("C39").Value = 4.1639999999999997
$ws.Range("D39").Value = 8.2929999999999993
$ws.Range("E39").Value = 13.055999999999999
$ws.Range("F39").Value = 16.724
$ws.Range("G39").Value = 20.513000000000002
$ws.Range("H39").Value = 24.638000000000002
$ws.Range("I39").Value = 29.126000000000001
$ws.Range("J39").Value = 33.177
$ws.Range("K39").Value = 37.290999999999997
$ws.Range("L39").Value = 41.436
$ws.Range("C40").Value = 7.1340000000000003
$ws.Range("D40").Value = 14.225
$ws.Range("E40").Value = 21.355
$ws.Range("F40").Value = 28.428999999999998
$ws.Range("G40").Value = 35.511000000000003
$ws.Range("H40").Value = 42.384
$ws.Range("I40").Value = 49.786000000000001
$ws.Range("J40").Value = 56.872999999999998
$ws.Range("K40").Value = 63.954000000000001
$ws.Range("L40").Value = 72.415999999999997
$ws.Range("C41").Value = 3.2360000000000002
$ws.Range("D41").Value = 6.415
$ws.Range("E41").Value = 9.7370000000000001
$ws.Range("F41").Value = 12.846
$ws.Range("G41").Value = 16.327999999999999
$ws.Range("H41").Value = 19.623999999999999
$ws.Range("I41").Value = 22.581
$ws.Range("J41").Value = 26.835999999999999
$ws.Range("K41").Value = 19.212
$ws.Range("L41").Value = 33.145000000000003
$ws.Range("C42").Value = 39.6
$ws.Range("D42").Value = 79.2
$ws.Range("E42").Value = 118.8
$ws.Range("F42").Value = 158.4
$ws.Range("G42").Value = 198
$ws.Range("H42").Value = 237.6
$ws.Range("I42").Value = 277.2
$ws.Range("J42").Value = 316.8
$ws.Range("K42").Value = 356.4
$ws.Range("L42").Value = 396
$ws.Range("C43").Value = 78
$ws.Range("D43").Value = 156
$ws.Range("E43").Value = 234
$ws.Range("F43").Value = 312
$ws.Range("G43").Value = 390
$ws.Range("H43").Value = 468
$ws.Range("I43").Value = 546
$ws.Range("J43").Value = 624
$ws.Range("K43").Value = 702
$ws.Range("L43").Value = 780
$ws.Range("C44").Value = 30
$ws.Range("D44").Value = 60
$ws.Range("E44").Value = 90
$ws.Range("F44").Value = 120
$ws.Range("G44").Value = 150
$ws.Range("H44").Value = 180
$ws.Range("I44").Value = 210
$ws.Range("J44").Value = 240
$ws.Range("K44").Value = 270
$ws.Range("L44").Value = 300
$ws.Range("E50").Value = 1

# F50 is rescaled too (interleaved with the G50 column removal in the diff)
$ws.Range("F50").Value = 1

# --- Row 48: drop G/H/I/L/M/N entirely; J/K/O keep their style but go blank ---
$ws.Range("G48").ClearContents()
$ws.Range("H48").ClearContents()
$ws.Range("I48").ClearContents()
$ws.Range("L48").ClearContents()
$ws.Range("M48").ClearContents()
$ws.Range("N48").ClearContents()
$ws.Range("J48").ClearContents()
$ws.Range("K48").ClearContents()
$ws.Range("O48").ClearContents()

# --- Row 49: drop G49:O49 entirely ---
$ws.Range("G49").ClearContents()
$ws.Range("H49").ClearContents()
$ws.Range("I49").ClearContents()
$ws.Range("J49").ClearContents()
$ws.Range("K49").ClearContents()
$ws.Range("L49").ClearContents()
$ws.Range("M49").ClearContents()
$ws.Range("N49").ClearContents()
$ws.Range("O49").ClearContents()

# --- Row 50: drop G50:O50 entirely ---
$ws.Range("G50").ClearContents()
$ws.Range("H50").ClearContents()
$ws.Range("I50").ClearContents()
$ws.Range("J50").ClearContents()
$ws.Range("K50").ClearContents()
$ws.Range("L50").ClearContents()
$ws.Range("M50").ClearContents()
$ws.Range("N50").ClearContents()
$ws.Range("O50").ClearContents()

# --- Update the active selection ---
$ws.Range("L23").Select()
